$d = $word.ActiveDocument

# Every paragraph in the document body gets an explicit
# pageBreakBefore = False setting in its paragraph properties.
$d.Paragraphs.PageBreakBefore = 0

# The same explicit setting is added to the built-in heading / title
# paragraph styles used by the document.
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $s = $d.Styles.Item($name)
    $s.ParagraphFormat.PageBreakBefore = 0
}

Write-Output "pageBreakBefore applied to paragraphs and styles"
